# Update Airbnb Price Model_CapStone.pptx
# 1. Refresh the cached "datetimeFigureOut" date placeholder text
#    (4/7/2020 -> 9/10/2020) everywhere it is used: the slide master,
#    every slide layout, and the notes master.
# 2. Fix the misspelled CONCLUSION slide title.

$p = $ppt.ActivePresentation

$newDate = "9/10/2020"

# --- Slide Master ---------------------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# --- Every slide layout ----------------------------------------------
$layouts = $master.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    $layout = $layouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Notes master ------------------------------------------------------
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# --- Slide 20: fix "COCLUSION" -> "CONCLUSION" title -------------------
$slide20 = $p.Slides.Item(20)
for ($i = 1; $i -le $slide20.Shapes.Count; $i++) {
    $shp = $slide20.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "COCLUSION") {
            $shp.TextFrame.TextRange.Text = "CONCLUSION"
        }
    }
}
